$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "287.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.62%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.12%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.058"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.62%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06659"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.22%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.364"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.52%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.393"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.376"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.87%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9428"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.45%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1557"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.17%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06627"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.08%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07636"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.27%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02954"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.14%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.08998"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.04%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001575"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.63%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04485"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.40%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0006450"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.95%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006317"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "5.06%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.445"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.18%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.263"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.58%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.27%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.95%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.062"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "4.25%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.22%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001181"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.43%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004490"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.98%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.98%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001617"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-2.19%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04211"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.48%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006746"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.64%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-11.03%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002019"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.84%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01235"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "12.07%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005688"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.39%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "25.93%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-29.36%"
